$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.916.80'
$ws.Range("E2").Value = '  -0.99%  '
$ws.Range("D3").Value = '2.340.32'
$ws.Range("E3").Value = '  +0.68%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '302.85'
$ws.Range("E5").Value = '  +0.13%  '
$ws.Range("D6").Value = '94.62'
$ws.Range("E6").Value = '  -3.59%  '
$ws.Range("D7").Value = '0.503'
$ws.Range("E7").Value = '  -1.12%  '
$ws.Range("E8").Value = '  -0.05%  '
$ws.Range("D9").Value = '0.496'
$ws.Range("E9").Value = '  -1.63%  '
$ws.Range("D10").Value = '34.14'
$ws.Range("E10").Value = '  -4.24%  '
$ws.Range("D11").Value = '0.0783'
$ws.Range("E11").Value = '  -1.80%  '
$ws.Range("D12").Value = '''18.70'
$ws.Range("E12").Value = '  -4.98%  '
$ws.Range("E13").Value = '  +1.51%  '
$ws.Range("D14").Value = '6.75'
$ws.Range("E14").Value = '  -2.43%  '
$ws.Range("D15").Value = '2.702.35'
$ws.Range("E15").Value = '  +0.47%  '
$ws.Range("D16").Value = '2.332.34'
$ws.Range("E16").Value = '  +0.39%  '
$ws.Range("D17").Value = '0.796'
$ws.Range("E17").Value = '  +0.31%  '
$ws.Range("D18").Value = '42.783.75'
$ws.Range("E18").Value = '  -1.09%  '
$ws.Range("D19").Value = '12.11'
$ws.Range("E19").Value = '  -4.92%  '
$ws.Range("D20").Value = '6.19'
$ws.Range("E20").Value = '  +1.76%  '
$ws.Range("E21").Value = '  -1.28%  '
$ws.Range("D22").Value = '67.89'
$ws.Range("E22").Value = '  -0.25%  '
$ws.Range("D23").Value = '235.72'
$ws.Range("E23").Value = '  -0.61%  '
$ws.Range("D24").Value = '2.21'
$ws.Range("E24").Value = '  -2.09%  '
$ws.Range("E25").Value = '  +0.02%  '
$ws.Range("E26").Value = '  -2.09%  '
$ws.Range("D27").Value = '24.62'
$ws.Range("E27").Value = '  -1.94%  '
$ws.Range("D28").Value = '2.34'
$ws.Range("E28").Value = '  +13.15%  '
$ws.Range("D29").Value = '9.15'
$ws.Range("E29").Value = '  +0.06%  '
$ws.Range("D30").Value = '31.37'
$ws.Range("E30").Value = '  -5.55%  '
$ws.Range("D31").Value = '0.999'
$ws.Range("E31").Value = '  +0.07%  '
$ws.Range("D32").Value = '''5.00'
$ws.Range("E32").Value = '  -0.35%  '
$ws.Range("D33").Value = '0.0734'
$ws.Range("E33").Value = '  +3.87%  '
$ws.Range("D34").Value = '17.28'
$ws.Range("E34").Value = '  -3.54%  '
$ws.Range("D38").Value = '124.81'
$ws.Range("E38").Value = '  -24.04%  '
$ws.Range("E39").Value = '  -0.62%  '
$ws.Range("D40").Value = '2.77'
$ws.Range("E40").Value = '  -1.10%  '
$ws.Range("D41").Value = '''22.10'
$ws.Range("E41").Value = '  +15.45%  '
$ws.Range("E42").Value = '  -1.39%  '
$ws.Range("D43").Value = '1.938.24'
$ws.Range("E43").Value = '  -2.42%  '
$ws.Range("E44").Value = '  -0.10%  '
$ws.Range("D45").Value = '10.13'
$ws.Range("E45").Value = '  -4.80%  '
$ws.Range("D46").Value = '2.09'
$ws.Range("E46").Value = '  +0.58%  '
$ws.Range("E47").Value = '  -3.40%  '
$ws.Range("D48").Value = '2.567.54'
$ws.Range("E48").Value = '  +0.44%  '
$ws.Range("E49").Value = '  +0.26%  '
$ws.Range("D50").Value = '52.84'
$ws.Range("E50").Value = '  -2.10%  '
$ws.Range("D51").Value = '71.75'
$ws.Range("E51").Value = '  -1.37%  '

# Row reordering: WEMIXToken/RenderToken/ARBITRUM rotate positions (rows 35-37)
$ws.Range("B35").Value = 'RenderToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D35").Value = '4.37'
$ws.Range("E35").Value = '  -3.17%  '
$ws.Range("B36").Value = 'ARBITRUM'
$ws.Range("C36").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D36").Value = '1.83'
$ws.Range("E36").Value = '  +2.35%  '
$ws.Range("B37").Value = 'WEMIXToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D37").Value = '2.33'
$ws.Range("E37").Value = '  -1.18%  '
